# The deck ships two theme parts: theme1.xml ("Office Theme", wired to the
# notes master) and theme2.xml ("Integral", wired to the slide master / the
# presentation's actual Design). The authored edit swaps the two themes'
# contents, so the slide master ends up on the stock "Office Theme" palette
# (its name/fontScheme/fmtScheme were already identical to theme1's, the
# only real payload is the 12-colour scheme).
#
# Recolour the active colour scheme (SlideMaster.ColorScheme, backed by
# ppt/theme/theme2.xml) to the "Office Theme" palette, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# PowerPoint's ColorScheme.Colors(n).RGB takes a COM BGR-packed long
# (R + G*256 + B*65536), matching the srgbClr hex RRGGBB below.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

$scheme.Colors(1).RGB  = 0        # dk1      000000
$scheme.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388  # dk2      44546A
$scheme.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501  # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407    # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308 # accent5  4472C4
$scheme.Colors(10).RGB = 4697456  # accent6  70AD47
$scheme.Colors(11).RGB = 12673797 # hlink    0563C1
$scheme.Colors(12).RGB = 7491477  # folHlink 954F72
